$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.344.74"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "1.869.74"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "'235.79"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "'0.4684"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("D9").Value = "'0.06546"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "'21.42"
$ws.Range("E10").Value = "  +6.37%  "
$ws.Range("D11").Value = "'0.07887"
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").Value = "'98.20"
$ws.Range("E12").Value = "  +2.45%  "
$ws.Range("D13").Value = "1.878.36"
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("D14").Value = "'5.111"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").Value = "'0.6767"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").Value = "'279.06"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "30.347.50"
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("E19").Value = "  +1.98%  "
$ws.Range("D20").Value = "'5.478"
$ws.Range("E20").Value = "  +2.54%  "
$ws.Range("D21").Value = "2.121.72"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "'0.000007312"
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").Value = "'6.157"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "'165.37"
$ws.Range("D26").Value = "'9.158"
$ws.Range("E26").Value = "  -1.35%  "
$ws.Range("D27").Value = "'19.18"
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("D28").Value = "'1.934"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").Value = "'1.382"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("D30").Value = "'0.09690"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").Value = "'4.397"
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("D32").Value = "'1.475"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").Value = "'4.104"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").Value = "'0.04708"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").Value = "'1.127"
$ws.Range("E35").Value = "  +4.34%  "
$ws.Range("D36").Value = "'0.7067"
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").Value = "'0.01861"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").Value = "'6.329"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").Value = "'2.538"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("D41").Value = "'74.50"
$ws.Range("E41").Value = "  +4.78%  "
$ws.Range("D42").Value = "'1.956"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("D43").Value = "'0.8507"
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("D44").Value = "'0.4189"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "'103.86"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("D47").Value = "'7.221"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").Value = "'9.265"
$ws.Range("E48").Value = "  +2.24%  "
$ws.Range("D49").Value = "'938.20"
$ws.Range("E49").Value = "  -3.94%  "
$ws.Range("D50").Value = "'34.24"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("E51").Value = "  -1.22%  "
